# Auto-generated Excel COM-interop script to apply scheduled-runner price refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H..N) across the
# Leve-profit tracker tables on each job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 469.2857
$ws.Range("I6").Value = 206.36363
$ws.Range("J6").Value = 1433.3334
$ws.Range("K6").Value = 619.0908899999999
$ws.Range("L6").Value = 4300.0002
$ws.Range("M6").Value = -507.0908899999999
$ws.Range("N6").Value = -4524.0002
$ws.Range("H64").Value = 3137.1428
$ws.Range("I64").Value = 3320
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3320
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -3072
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3137.1428
$ws.Range("I67").Value = 3320
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3320
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2462
$ws.Range("N67").Value = -4716
$ws.Range("H86").Value = 2984
$ws.Range("I86").Value = 2666.6667
$ws.Range("J86").Value = 3301.3333
$ws.Range("K86").Value = 2666.6667
$ws.Range("L86").Value = 3301.3333
$ws.Range("M86").Value = -1543.6667
$ws.Range("N86").Value = -5547.3333
$ws.Range("H89").Value = 2984
$ws.Range("I89").Value = 2666.6667
$ws.Range("J89").Value = 3301.3333
$ws.Range("K89").Value = 13333.3335
$ws.Range("L89").Value = 16506.6665
$ws.Range("M89").Value = -7717.333500000001
$ws.Range("N89").Value = -27738.6665
$ws.Range("H132").Value = 2533342
$ws.Range("I132").Value = 2668138.8
$ws.Range("J132").Value = 5899.5
$ws.Range("K132").Value = 8004416.399999999
$ws.Range("L132").Value = 17698.5
$ws.Range("M132").Value = -8001886.399999999
$ws.Range("N132").Value = -22758.5
$ws.Range("H133").Value = 29738.182
$ws.Range("J133").Value = 29738.182
$ws.Range("L133").Value = 29738.182
$ws.Range("N133").Value = -39858.182
$ws.Range("H138").Value = 1938.97
$ws.Range("I138").Value = 620.2820400000001
$ws.Range("J138").Value = 2782.0657
$ws.Range("K138").Value = 1860.84612
$ws.Range("L138").Value = 8346.197100000001
$ws.Range("M138").Value = 3279.15388
$ws.Range("N138").Value = -18626.1971
$ws.Range("H141").Value = 384841.94
$ws.Range("I141").Value = 1483.3334
$ws.Range("J141").Value = 2224963.2
$ws.Range("K141").Value = 4450.0002
$ws.Range("L141").Value = 6674889.600000001
$ws.Range("M141").Value = 729.9997999999996
$ws.Range("N141").Value = -6685249.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 72319
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10346
$ws.Range("H32").Value = 1042.46
$ws.Range("I32").Value = 941.5326
$ws.Range("J32").Value = 2203.125
$ws.Range("K32").Value = 941.5326
$ws.Range("L32").Value = 2203.125
$ws.Range("M32").Value = -654.5326
$ws.Range("N32").Value = -2777.125
$ws.Range("H63").Value = 6999.75
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6999.75
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 6999.75
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -8371.75
$ws.Range("H66").Value = 6999.75
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6999.75
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 34998.75
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -41862.75
$ws.Range("H74").Value = 554.1539
$ws.Range("I74").Value = 519.2432
$ws.Range("J74").Value = 1200
$ws.Range("K74").Value = 519.2432
$ws.Range("L74").Value = 1200
$ws.Range("M74").Value = 354.7568
$ws.Range("N74").Value = -2948
$ws.Range("H77").Value = 554.1539
$ws.Range("I77").Value = 519.2432
$ws.Range("J77").Value = 1200
$ws.Range("K77").Value = 2596.216
$ws.Range("L77").Value = 6000
$ws.Range("M77").Value = 1771.784
$ws.Range("N77").Value = -14736
$ws.Range("H97").Value = 858
$ws.Range("I97").Value = 862.125
$ws.Range("J97").Value = 849.75
$ws.Range("K97").Value = 862.125
$ws.Range("L97").Value = 849.75
$ws.Range("M97").Value = -366.125
$ws.Range("N97").Value = -1841.75
$ws.Range("H132").Value = 1819.3889
$ws.Range("I132").Value = 1322.2972
$ws.Range("J132").Value = 2901.2942
$ws.Range("K132").Value = 3966.8916
$ws.Range("L132").Value = 8703.882599999999
$ws.Range("M132").Value = -1436.8916
$ws.Range("N132").Value = -13763.8826

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 832.96875
$ws.Range("I94").Value = 749.7727
$ws.Range("J94").Value = 1016
$ws.Range("K94").Value = 749.7727
$ws.Range("L94").Value = 1016
$ws.Range("M94").Value = -298.7727
$ws.Range("N94").Value = -1918
$ws.Range("H126").Value = 25125
$ws.Range("H134").Value = 2503.7878
$ws.Range("I134").Value = 1723.3636
$ws.Range("J134").Value = 4064.6365
$ws.Range("K134").Value = 5170.0908
$ws.Range("L134").Value = 12193.9095
$ws.Range("M134").Value = -2635.0908
$ws.Range("N134").Value = -17263.9095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 537.8570999999999
$ws.Range("I19").Value = 294.16666
$ws.Range("J19").Value = 2000
$ws.Range("K19").Value = 294.16666
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = -124.16666
$ws.Range("N19").Value = -2340
$ws.Range("H24").Value = 537.8570999999999
$ws.Range("I24").Value = 294.16666
$ws.Range("J24").Value = 2000
$ws.Range("K24").Value = 294.16666
$ws.Range("L24").Value = 2000
$ws.Range("M24").Value = -124.16666
$ws.Range("N24").Value = -2340
$ws.Range("H58").Value = 8476596
$ws.Range("I58").Value = 958.4681
$ws.Range("J58").Value = 41672844
$ws.Range("K58").Value = 958.4681
$ws.Range("L58").Value = 41672844
$ws.Range("M58").Value = -755.4681
$ws.Range("N58").Value = -41673250
$ws.Range("H134").Value = 1389.0227
$ws.Range("I134").Value = 715.46155
$ws.Range("K134").Value = 2146.38465
$ws.Range("M134").Value = 388.61535
$ws.Range("H136").Value = 8476596
$ws.Range("I136").Value = 958.4681
$ws.Range("J136").Value = 41672844
$ws.Range("K136").Value = 2875.4043
$ws.Range("L136").Value = 125018532
$ws.Range("M136").Value = -325.4043000000001
$ws.Range("N136").Value = -125023632

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H118").Value = 1914.375
$ws.Range("I118").Value = 800
$ws.Range("J118").Value = 2583
$ws.Range("K118").Value = 2400
$ws.Range("L118").Value = 7749
$ws.Range("M118").Value = -1157
$ws.Range("N118").Value = -10235

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3688.7693
$ws.Range("I97").Value = 2541.3635
$ws.Range("J97").Value = 9999.5
$ws.Range("K97").Value = 2541.3635
$ws.Range("L97").Value = 9999.5
$ws.Range("M97").Value = -2045.3635
$ws.Range("N97").Value = -10991.5
$ws.Range("H113").Value = 3488.8333
$ws.Range("I113").Value = 2844.3333
$ws.Range("J113").Value = 4133.3335
$ws.Range("K113").Value = 2844.3333
$ws.Range("L113").Value = 4133.3335
$ws.Range("M113").Value = -674.3332999999998
$ws.Range("N113").Value = -8473.333500000001
$ws.Range("H122").Value = 4053.9375
$ws.Range("I122").Value = 2875.7334
$ws.Range("J122").Value = 5093.5293
$ws.Range("K122").Value = 8627.200199999999
$ws.Range("L122").Value = 15280.5879
$ws.Range("M122").Value = -6177.200199999999
$ws.Range("N122").Value = -20180.5879
$ws.Range("H132").Value = 1858.4717
$ws.Range("I132").Value = 1462.4
$ws.Range("J132").Value = 3077.1538
$ws.Range("K132").Value = 4387.200000000001
$ws.Range("L132").Value = 9231.4614
$ws.Range("M132").Value = -1857.200000000001
$ws.Range("N132").Value = -14291.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2892.9333
$ws.Range("I93").Value = 2499.1667
$ws.Range("K93").Value = 2499.1667
$ws.Range("M93").Value = -1251.1667
$ws.Range("H122").Value = 2755.2258
$ws.Range("I122").Value = 2287.4783
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 6862.4349
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -4412.4349
$ws.Range("N122").Value = -17200
$ws.Range("H132").Value = 1910.1538
$ws.Range("I132").Value = 1265.5
$ws.Range("J132").Value = 4857.143
$ws.Range("K132").Value = 3796.5
$ws.Range("L132").Value = 14571.429
$ws.Range("M132").Value = -1266.5
$ws.Range("N132").Value = -19631.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 246913.72
$ws.Range("I132").Value = 296764.12
$ws.Range("J132").Value = 58590
$ws.Range("K132").Value = 890292.36
$ws.Range("L132").Value = 175770
$ws.Range("M132").Value = -887762.36
$ws.Range("N132").Value = -180830
